$d = $word.ActiveDocument

# --- 1. Title paragraph: drop the first-line indent and swap the text ---
# " БИБЛИОТЕКИ ГОРОДА" -> "УЧЕТ АВТОМОБИЛЕЙ" (the run's xml:space="preserve"
# goes away on its own because the new text has no leading/trailing space)
$found = $d.Content.Find.Execute(" БИБЛИОТЕКИ ГОРОДА", $true, $false, $false,
    $false, $false, $true, 1, $false, "УЧЕТ АВТОМОБИЛЕЙ", 2)

# Locate that paragraph again (Find collapses the range) and zero its
# first-line indent, which drops the <w:ind> from the paragraph.
$p = $d.Content.Find.Parent.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*АВТОМОБИЛЕЙ*") {
        $cand.Format.FirstLineIndent = 0
        break
    }
}

# --- 2. Move the stray "_GoBack" bookmark up onto the empty paragraph
#        that follows the title (Word keeps bookmark names unique, so
#        re-adding it elsewhere relocates it and renumbers the other
#        bookmark ids the same way real Word does) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*АВТОМОБИЛЕЙ*") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}
$d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
